{"js": "// The edit removes the two trailing empty paragraphs that follow the\n// final line of dialogue (\"Actually, now that I think about it, she\n// probably practices on Saturday as well.\") at the end of the document\n// body, right before the section break.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst anchorText =\n  \"Actually, now that I think about it, she probably practices on Saturday as well.\";\n\n// Find the anchor paragraph (the last paragraph that actually has text).\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\n// Remove every empty paragraph that trails the anchor paragraph (walking\n// backwards so the indices of not-yet-deleted items stay valid).\nif (anchorIndex !== -1) {\n  for (let i = items.length - 1; i > anchorIndex; i--) {\n    if (items[i].text === \"\") {\n      items[i].delete();\n    }\n  }\n} else {\n  // Fallback: if the anchor text couldn't be located (e.g. formatting\n  // differences), just drop trailing empty paragraphs from the end of\n  // the body.\n  for (let i = items.length - 1; i >= 0 && items[i].text === \"\"; i--) {\n    items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# The edit removes the two trailing empty paragraphs that follow the\n# final line of dialogue (\"Actually, now that I think about it, she\n# probably practices on Saturday as well.\") at the end of the document\n# body, right before the section break.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Actually, now that I think about it, she probably practices on Saturday as well.\"\n\n# Locate the anchor paragraph so we know where the trailing empty\n# paragraphs start.\n$anchorIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13)\n    if ($text -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\n# Delete every empty trailing paragraph after the anchor, walking\n# backwards from the end of the document so indices stay valid as\n# paragraphs are removed.\nfor ($i = $d.Paragraphs.Count; $i -gt $anchorIndex; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -eq [char]13) {\n        $p.Range.Delete()\n    }\n}\n"}
